$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$origStyleD = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.673.16"
$ws.Range("D2").Style = $origStyleD
$ws.Range("E2").Value = "  -1.05%  "

# Row 3
$origStyleD = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.787.39"
$ws.Range("D3").Style = $origStyleD
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$origStyleD = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.94"
$ws.Range("D5").Style = $origStyleD
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$origStyleD = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.15"
$ws.Range("D6").Style = $origStyleD
$ws.Range("E6").Value = "  +0.35%  "

# Row 7
$origStyleD = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.772.19"
$ws.Range("D7").Style = $origStyleD
$ws.Range("E7").Value = "  +0.68%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("E9").Value = "  +0.11%  "

# Row 10
$ws.Range("E10").Value = "  -0.38%  "

# Row 11
$origStyleD = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("D11").Style = $origStyleD
$ws.Range("E11").Value = "  -1.90%  "

# Row 12
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("E13").Value = "  -3.12%  "

# Row 14
$origStyleD = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.97"
$ws.Range("D14").Style = $origStyleD
$ws.Range("E14").Value = "  -0.82%  "

# Row 15
$origStyleD = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.422.83"
$ws.Range("D15").Style = $origStyleD
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$origStyleD = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.785.82"
$ws.Range("D16").Style = $origStyleD
$ws.Range("E16").Value = "  +1.40%  "

# Row 17
$origStyleD = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.53"
$ws.Range("D17").Style = $origStyleD
$ws.Range("E17").Value = "  +3.66%  "

# Row 18
$origStyleD = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.638.91"
$ws.Range("D18").Style = $origStyleD
$ws.Range("E18").Value = "  -1.06%  "

# Row 19
$origStyleD = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("D19").Style = $origStyleD
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("E20").Value = "  -0.19%  "

# Row 21
$origStyleD = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.01"
$ws.Range("D21").Style = $origStyleD
$ws.Range("E21").Value = "  -6.16%  "

# Row 22
$origStyleD = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.20"
$ws.Range("D22").Style = $origStyleD
$ws.Range("E22").Value = "  -1.32%  "

# Row 23
$origStyleD = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").Style = $origStyleD
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
$origStyleD = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("D24").Style = $origStyleD
$ws.Range("E24").Value = "  +3.69%  "

# Row 25
$ws.Range("E25").Value = "  -0.75%  "

# Row 26
$origStyleD = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.00"
$ws.Range("D26").Style = $origStyleD
$ws.Range("E26").Value = "  +0.58%  "

# Row 27
$ws.Range("E27").Value = "  -3.13%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$origStyleD = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = $origStyleD
$ws.Range("E28").Value = "  -0.83%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$origStyleD = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = $origStyleD
$ws.Range("E29").Value = "  +0.11%  "

# Row 30
$origStyleD = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.933.08"
$ws.Range("D30").Style = $origStyleD
$ws.Range("E30").Value = "  +0.98%  "

# Row 31
$ws.Range("E31").Value = "  +0.41%  "

# Row 32
$ws.Range("E32").Value = "  +2.99%  "

# Row 33
$origStyleD = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.19"
$ws.Range("D33").Style = $origStyleD
$ws.Range("E33").Value = "  -1.80%  "

# Row 34
$origStyleD = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.62"
$ws.Range("D34").Style = $origStyleD
$ws.Range("E34").Value = "  -1.12%  "

# Row 35
$origStyleD = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = $origStyleD

# Row 36
$origStyleD = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.10"
$ws.Range("D36").Style = $origStyleD
$ws.Range("E36").Value = "  -1.05%  "

# Row 37
$origStyleD = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0999"
$ws.Range("D37").Style = $origStyleD
$ws.Range("E37").Value = "  -0.77%  "

# Row 38
$ws.Range("E38").Value = "  -3.03%  "

# Row 39
$ws.Range("E39").Value = "  -0.17%  "

# Row 40
$origStyleD = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.994"
$ws.Range("D40").Style = $origStyleD
$ws.Range("E40").Value = "  -0.29%  "

# Row 41
$origStyleD = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("D41").Style = $origStyleD
$ws.Range("E41").Value = "  -0.41%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$origStyleD = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.82"
$ws.Range("D44").Style = $origStyleD
$ws.Range("E44").Value = "  +5.72%  "

# Row 45
$origStyleD = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.15"
$ws.Range("D45").Style = $origStyleD
$ws.Range("E45").Value = "  +3.15%  "

# Row 46
$ws.Range("E46").Value = "  -1.07%  "

# Row 47
$origStyleD = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.02"
$ws.Range("D47").Style = $origStyleD
$ws.Range("E47").Value = "  +4.32%  "

# Row 48
$ws.Range("E48").Value = "  -1.92%  "

# Row 49
$origStyleD = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "393.77"
$ws.Range("D49").Style = $origStyleD
$ws.Range("E49").Value = "  +0.92%  "

# Row 50
$origStyleD = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.76"
$ws.Range("D50").Style = $origStyleD
$ws.Range("E50").Value = "  +6.44%  "

# Row 51
$ws.Range("E51").Value = "  -5.39%  "
